# Quarterly-separated product sheet: roll the quarter window forward by one
# quarter. Every quarterly data row (columns E:N) drops its oldest quarter
# (column E) and gains a new quarter at the end (column N); the shared
# "quarter label" header rows (8/28/48/61/81) pick up the new label
# automatically because they reference the same shared-string slot whose
# text is rotated the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E","F","G","H","I","J","K","L","M","N")

# New value to place in column N (the newest quarter, "فصل چهارم منتهی به
# 1401/12") for every affected data row. "-" marks the textual placeholder
# used throughout the sheet for "no data".
$newN = @{
    11 = 706665
    12 = 50940
    13 = 1292885
    14 = 2050490
    16 = 0
    17 = 2233
    18 = 0
    19 = 2233
    21 = 0
    22 = "-"
    23 = 0
    24 = 2052723
    31 = 3587536
    32 = 282318
    33 = 320076
    34 = 4189930
    36 = 0
    37 = 375691
    38 = 0
    39 = 375691
    41 = 0
    42 = "-"
    43 = 0
    44 = 4565621
    51 = 5076714
    52 = 5542167
    53 = 247567
    55 = "-"
    56 = 168244962
    57 = "-"
    64 = -1993509
    65 = -283429
    66 = -265901
    67 = -2542839
    69 = 0
    70 = -89079
    71 = 0
    72 = -89079
    74 = 0
    75 = "-"
    76 = 0
    77 = -2631918
    84 = 1594027
    85 = -1111
    86 = 54175
    87 = 1647091
    89 = 0
    90 = 286612
    91 = 0
    92 = 286612
    94 = 0
    95 = 1933703
}

foreach ($row in $newN.Keys) {
    $vals = @()
    foreach ($c in $cols) {
        $vals += $ws.Range($c + $row).Value2
    }
    for ($i = 0; $i -lt 9; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i + 1]
    }
    $ws.Range("N" + $row).Value = $newN[$row]
}

# The header rows that carry the quarter-label strings (row 8 and its
# repeats at 28/48/61/81) each hold the same ten quarter labels, so rotate
# the label text in every one of them the same way the data rotated.
$labelCols = @("E","F","G","H","I","J","K","L","M","N")
$labelRows = @(8, 28, 48, 61, 81)
foreach ($lrow in $labelRows) {
    $labels = @()
    foreach ($c in $labelCols) {
        $labels += $ws.Range($c + $lrow).Value2
    }
    for ($i = 0; $i -lt 9; $i++) {
        $ws.Range($labelCols[$i] + $lrow).Value = $labels[$i + 1]
    }
    $ws.Range("N" + $lrow).Value = "فصل چهارم منتهی به 1401/12"
}
